# Auto-generated Excel COM-interop script applying the Ravana_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 135
$ws.Range("H135").Value = 3819.7778
$ws.Range("I135").Value = 4506
$ws.Range("K135").Value = 40554
$ws.Range("M135").Value = -38019

# Row 137
$ws.Range("H137").Value = 1889.0952
$ws.Range("I137").Value = 1282.8462
$ws.Range("J137").Value = 2874.25
$ws.Range("K137").Value = 3848.5386
$ws.Range("L137").Value = 8622.75
$ws.Range("M137").Value = -1298.5386
$ws.Range("N137").Value = -13722.75

# Row 138
$ws.Range("H138").Value = 7166.3447
$ws.Range("I138").Value = 3296.4
$ws.Range("K138").Value = 9889.200000000001
$ws.Range("M138").Value = -4749.200000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 22
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# Row 61
$ws.Range("H61").Value = 3817.5715
$ws.Range("I61").Value = 3745.2
$ws.Range("J61").Value = 3998.5
$ws.Range("K61").Value = 3745.2
$ws.Range("L61").Value = 3998.5
$ws.Range("M61").Value = -3533.2
$ws.Range("N61").Value = -4422.5

# Row 74
$ws.Range("H74").Value = 1721
$ws.Range("I74").Value = 1007.8333
$ws.Range("K74").Value = 1007.8333
$ws.Range("M74").Value = -133.8333

# Row 77
$ws.Range("H77").Value = 1721
$ws.Range("I77").Value = 1007.8333
$ws.Range("K77").Value = 5039.1665
$ws.Range("M77").Value = -671.1665000000003

# Row 97
$ws.Range("H97").Value = 661.5
$ws.Range("I97").Value = 712.7778
$ws.Range("K97").Value = 712.7778
$ws.Range("M97").Value = -216.7778

# Row 110
$ws.Range("H110").Value = 3562
$ws.Range("I110").Value = 3606.0833
$ws.Range("J110").Value = 3297.5
$ws.Range("K110").Value = 3606.0833
$ws.Range("L110").Value = 3297.5
$ws.Range("M110").Value = -1561.0833
$ws.Range("N110").Value = -7387.5

# Row 132
$ws.Range("H132").Value = 3698.389
$ws.Range("I132").Value = 3175.2222
$ws.Range("K132").Value = 9525.6666
$ws.Range("M132").Value = -6995.6666

# Row 136
$ws.Range("H136").Value = 3817.5715
$ws.Range("I136").Value = 3745.2
$ws.Range("J136").Value = 3998.5
$ws.Range("K136").Value = 11235.6
$ws.Range("L136").Value = 11995.5
$ws.Range("M136").Value = -8685.599999999999
$ws.Range("N136").Value = -17095.5

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498

# Row 107
$ws.Range("H107").Value = 2156.5
$ws.Range("I107").Value = 2156.5
$ws.Range("K107").Value = 2156.5
$ws.Range("M107").Value = -236.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4044.4614
$ws.Range("I31").Value = 1651.3529
$ws.Range("K31").Value = 1651.3529
$ws.Range("M31").Value = -1356.3529

# Row 34
$ws.Range("H34").Value = 4044.4614
$ws.Range("I34").Value = 1651.3529
$ws.Range("K34").Value = 1651.3529
$ws.Range("M34").Value = -1449.3529

# Row 82
$ws.Range("H82").Value = 25000
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25722

# Row 85
$ws.Range("H85").Value = 25000
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27496

# Row 132
$ws.Range("H132").Value = 4735.143
$ws.Range("I132").Value = 4200.143
$ws.Range("J132").Value = 5270.143
$ws.Range("K132").Value = 12600.429
$ws.Range("L132").Value = 15810.429
$ws.Range("M132").Value = -10070.429
$ws.Range("N132").Value = -20870.429

# Row 133
$ws.Range("H133").Value = 124597.5
$ws.Range("J133").Value = 124597.5
$ws.Range("L133").Value = 124597.5
$ws.Range("N133").Value = -129657.5

$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# Row 68
$ws.Range("H68").Value = 1200
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 1200
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 117
$ws.Range("H117").Value = 1315
$ws.Range("J117").Value = 570
$ws.Range("L117").Value = 1710
$ws.Range("N117").Value = -8594

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("M2").Value = 112

# Row 102
$ws.Range("H102").Value = 2179.4
$ws.Range("I102").Value = 2179.4
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2179.4
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -557.4000000000001
$ws.Range("N102").ClearContents()

# Row 132
$ws.Range("H132").Value = 5333.769
$ws.Range("I132").Value = 5218.2856
$ws.Range("J132").Value = 5468.5
$ws.Range("K132").Value = 15654.8568
$ws.Range("L132").Value = 16405.5
$ws.Range("M132").Value = -13124.8568
$ws.Range("N132").Value = -21465.5

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2773
$ws.Range("I46").Value = 2553
$ws.Range("J46").Value = 2993
$ws.Range("K46").Value = 2553
$ws.Range("L46").Value = 2993
$ws.Range("M46").Value = -2365
$ws.Range("N46").Value = -3369

# Row 132
$ws.Range("H132").Value = 5996.8335
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5996.8335
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17990.5005
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -23050.5005

# Row 136
$ws.Range("H136").Value = 22580.875
$ws.Range("J136").Value = 32367.25
$ws.Range("L136").Value = 97101.75
$ws.Range("N136").Value = -102201.75

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4499.6665
$ws.Range("J122").Value = 6249.5
$ws.Range("L122").Value = 18748.5
$ws.Range("N122").Value = -23648.5

# Row 132
$ws.Range("H132").Value = 3255.16
$ws.Range("I132").Value = 3008.7058
$ws.Range("J132").Value = 3778.875
$ws.Range("K132").Value = 9026.117400000001
$ws.Range("L132").Value = 11336.625
$ws.Range("M132").Value = -6496.117400000001
$ws.Range("N132").Value = -16396.625

# Row 136
$ws.Range("H136").Value = 11446.846
$ws.Range("I136").Value = 13115.3
$ws.Range("K136").Value = 39345.89999999999
$ws.Range("M136").Value = -36795.89999999999
